$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "octopus"

$ws.Range("C14").Value = 45008.770833333336
$ws.Range("D14").Value = 45008.791666666664

$ws.Range("C14:D14").NumberFormat = "m/d/yy h:mm"

$ws.Range("D14").Select()
$excel.ActiveWindow.ScrollRow = 3
